# gamedata-faf-waw.xlsx -- "Added Techs in the excel file and import"
#
# 1) Fix the "Wonders " sheet name (trailing space removed).
# 2) Add four new "Level N Tech" sheets at the end of the workbook, each
#    with a bold "Name"/"Description" header row and a single column of
#    technology names pulled from the Civ V tech tree.
# 3) Leave "Level 4 Tech" as the active/selected sheet (mirrors the
#    workbook's new activeTab), which also naturally drops tabSelected
#    from the previously-active "City-states" sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Rename "Wonders " -> "Wonders" ---------------------------------
$wonders = $wb.Worksheets.Item("Wonders ")
$wonders.Name = "Wonders"

# --- 2) Add the four tech sheets ---------------------------------------
$techLevels = @(
    @{
        Name = "Level 1 Tech"
        Items = @(
            "Animal Husbandry",
            "Agriculture",
            "Code of Laws",
            "Currency",
            "Horseback Riding",
            "Masonry",
            "Metalworking",
            "Navigation",
            "Philosophy",
            "Pottery",
            "Writing",
            "Navy"
        )
        Selection = "D18"
    },
    @{
        Name = "Level 2 Tech"
        Items = @(
            "Civil Service",
            "Chivalry",
            "Construction",
            "Democracy",
            "Engineering",
            "Irrigation",
            "Mathematics",
            "Monarchy",
            "Mysticism",
            "Printing Press",
            "Sailing",
            "Logistics",
            "Bureaucracy"
        )
        Selection = "A2:A14"
    },
    @{
        Name = "Level 3 Tech"
        Items = @(
            "Banking",
            "Biology",
            "Communism",
            "Ecology",
            "Gunpowder",
            "Metal Casting",
            "Military Science",
            "Railroad",
            "Steam Power",
            "Theology",
            "Education"
        )
        Selection = "B8"
    },
    @{
        Name = "Level 4 Tech"
        Items = @(
            "Atomic Theory",
            "Ballistics",
            "Combustion",
            "Computers",
            "Flight",
            "Mass Media",
            "Plastics",
            "Replacement Parts",
            "Plastics"
        )
        Selection = "M47"
    }
)

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheets = @()

foreach ($lvl in $techLevels) {
    $ws = $wb.Worksheets.Add($null, $afterSheet)
    $ws.Name = $lvl.Name

    $ws.Range("A1").Value = "Name"
    $ws.Range("B1").Value = "Description"
    $ws.Range("A1:B1").Font.Bold = $true

    $row = 2
    foreach ($item in $lvl.Items) {
        $ws.Cells.Item($row, 1).Value = $item
        $row = $row + 1
    }

    $ws.Columns.Item(1).AutoFit() | Out-Null

    $newSheets += $ws
    $afterSheet = $ws
}

# --- 3) Selections + final active sheet ---------------------------------
for ($i = 0; $i -lt $newSheets.Count; $i++) {
    $ws = $newSheets[$i]
    $lvl = $techLevels[$i]
    $ws.Activate()
    $ws.Range($lvl.Selection).Select() | Out-Null
}

# "Level 4 Tech" (the last sheet added) ends up active, matching the
# workbook's new activeTab and clearing tabSelected on "City-states".
